$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name and Link for reordered rows (38-41, 48-49)
$bcUpdates = @{
    38 = @('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar')
    39 = @('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo')
    40 = @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp')
    41 = @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt')
    48 = @('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt')
    49 = @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near')
}
foreach ($r in $bcUpdates.Keys) {
    $pair = $bcUpdates[$r]
    $ws.Range("B$r").Value = $pair[0]
    $ws.Range("C$r").Value = $pair[1]
}

# Update Price column (D) - force text storage via leading apostrophe
# so numeric-looking strings are not reinterpreted as numbers
$dUpdates = @{
    2 = '22.388.18'
    3 = '1.560.71'
    4 = '1.001'
    6 = '284.23'
    7 = '0.3608'
    8 = '48.20'
    9 = '0.3321'
    10 = '1.123'
    11 = '0.07371'
    12 = '1.001'
    13 = '20.69'
    14 = '5.907'
    15 = '6.866'
    16 = '1.561.15'
    17 = '0.00001100'
    18 = '87.82'
    19 = '0.06688'
    20 = '1.001'
    21 = '6.325'
    22 = '16.03'
    23 = '11.96'
    24 = '22.387.62'
    25 = '2.416'
    26 = '2.535'
    28 = '19.33'
    29 = '4.988'
    30 = '122.70'
    31 = '1.736.54'
    32 = '1.053'
    33 = '6.104'
    34 = '1.996'
    35 = '9.767'
    36 = '0.08276'
    37 = '0.02382'
    38 = '0.06365'
    39 = '0.2201'
    40 = '5.290'
    41 = '1.269'
    42 = '11.08'
    43 = '0.6039'
    44 = '1.000'
    45 = '13.75'
    46 = '3.746'
    47 = '0.5733'
    48 = '124.20'
    49 = '1.997'
    50 = '1.206'
    51 = '0.07199'
}
foreach ($r in $dUpdates.Keys) {
    $ws.Range("D$r").Value = "'" + $dUpdates[$r]
}

# Update Volume(1h) column (E)
$eUpdates = @{
    2 = '  -0.21%  '
    3 = '  -0.78%  '
    4 = '  -0.20%  '
    5 = '  -0.14%  '
    6 = '  -2.46%  '
    7 = '  -3.33%  '
    8 = '  -3.34%  '
    9 = '  -2.16%  '
    10 = '  -1.07%  '
    11 = '  -2.44%  '
    12 = '  -0.18%  '
    13 = '  -3.34%  '
    14 = '  -1.44%  '
    15 = '  -0.90%  '
    16 = '  -0.87%  '
    17 = '  -2.07%  '
    18 = '  -3.50%  '
    19 = '  -0.78%  '
    20 = '  -0.16%  '
    21 = '  +0.84%  '
    22 = '  -2.16%  '
    23 = '  -1.41%  '
    24 = '  -0.25%  '
    25 = '  +3.05%  '
    26 = '  -2.83%  '
    27 = '  +0.74%  '
    28 = '  -3.85%  '
    29 = '  -0.33%  '
    30 = '  -2.42%  '
    31 = '  -0.84%  '
    32 = '  +0.23%  '
    33 = '  -0.56%  '
    34 = '  +0.77%  '
    35 = '  -0.26%  '
    36 = '  -1.43%  '
    37 = '  -3.22%  '
    38 = '  -2.29%  '
    39 = '  -3.82%  '
    40 = '  -3.07%  '
    41 = '  -7.89%  '
    42 = '  -1.86%  '
    43 = '  -2.91%  '
    44 = '  -0.16%  '
    45 = '  -1.21%  '
    46 = '  -1.68%  '
    47 = '  -1.24%  '
    48 = '  -4.24%  '
    49 = '  -3.74%  '
    50 = '  -0.89%  '
    51 = '  -1.62%  '
}
foreach ($r in $eUpdates.Keys) {
    $ws.Range("E$r").Value = $eUpdates[$r]
}